$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: trim the sentence ending "...by September 13, 2022, and
# report to jail on September 16, 2022, at 7:00 p.m." down to
# "...by September 13, 2022." -- i.e. remove the jail-reporting clause.
# ---------------------------------------------------------------------
$dateEnd = $d.Content
$dateEnd.Find.Execute("September 13, 2022", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterDate = $dateEnd.End

$sentenceEnd = $d.Content
$sentenceEnd.Find.Execute("p.m.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$periodEnd = $sentenceEnd.End

$toRemove1 = $d.Range($afterDate, $periodEnd)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"></w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Palatino Linotype" w:hAnsi="Palatino Linotype"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$toRemove1.InsertXML($xml1)

# ---------------------------------------------------------------------
# Change 2: delete the whole "Restitution." paragraph content (the
# boilerplate about paying $5,000 restitution to Justin Kudela),
# leaving the paragraph starting directly with "Fines and Costs."
# ---------------------------------------------------------------------
$restStart = $d.Content
$restStart.Find.Execute("Restitution.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startDel = $restStart.Start

$finesStart = $d.Content
$finesStart.Find.Execute("Fines and Costs.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endDel = $finesStart.Start

$toRemove2 = $d.Range($startDel, $endDel)
$toRemove2.Delete()
